# Update class and character
# Adds two new enums to the Sheet1 enum table:
#   SLG.Gender     (rows 29-31): NONE, MALE, FEMALE
#   SLG.ClassType  (rows 32-40): None, Armor, Backup, Cavalry, Covert,
#                                 Dragon, Flying, Mystical, QiAdept

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prime formatting for the new rows by cloning the existing
#     "SLG.WeaponRank" block (row 13, the template for a fresh enum
#     block: B/C/D carry the enum-name/flags/unique cells, G/I carry the
#     item name/value columns) before overwriting the values. ---
$ws.Range("B13").Copy($ws.Range("B29"))
$ws.Range("C13").Copy($ws.Range("C29"))
$ws.Range("D13").Copy($ws.Range("D29"))
$ws.Range("G13").Copy($ws.Range("G29"))
$ws.Range("I13").Copy($ws.Range("I29"))

$ws.Range("G13").Copy($ws.Range("G30"))
$ws.Range("I13").Copy($ws.Range("I30"))

$ws.Range("G13").Copy($ws.Range("G31"))
$ws.Range("I13").Copy($ws.Range("I31"))

$ws.Range("B13").Copy($ws.Range("B32"))
$ws.Range("C13").Copy($ws.Range("C32"))
$ws.Range("D13").Copy($ws.Range("D32"))
$ws.Range("G13").Copy($ws.Range("G32"))
$ws.Range("I13").Copy($ws.Range("I32"))

$ws.Range("G13").Copy($ws.Range("G33"))
$ws.Range("I13").Copy($ws.Range("I33"))
$ws.Range("G13").Copy($ws.Range("G34"))
$ws.Range("I13").Copy($ws.Range("I34"))
$ws.Range("G13").Copy($ws.Range("G35"))
$ws.Range("I13").Copy($ws.Range("I35"))
$ws.Range("G13").Copy($ws.Range("G36"))
$ws.Range("I13").Copy($ws.Range("I36"))
$ws.Range("G13").Copy($ws.Range("G37"))
$ws.Range("I13").Copy($ws.Range("I37"))
$ws.Range("G13").Copy($ws.Range("G38"))
$ws.Range("I13").Copy($ws.Range("I38"))
$ws.Range("G13").Copy($ws.Range("G39"))
$ws.Range("I13").Copy($ws.Range("I39"))
$ws.Range("G13").Copy($ws.Range("G40"))
$ws.Range("I13").Copy($ws.Range("I40"))

# --- SLG.Gender enum (rows 29-31) ---
$ws.Range("B29").Value = "SLG.Gender"
$ws.Range("C29").Value = $false
$ws.Range("D29").Value = $true
$ws.Range("G29").Value = "NONE"
$ws.Range("I29").Value = 0

$ws.Range("G30").Value = "MALE"
$ws.Range("I30").Value = 1

$ws.Range("G31").Value = "FEMALE"
$ws.Range("I31").Value = 2

# --- SLG.ClassType enum (rows 32-40) ---
# Enum header + values entered first without the "None" member ...
$ws.Range("B32").Value = "SLG.ClassType"
$ws.Range("C32").Value = $false
$ws.Range("D32").Value = $true

$ws.Range("G33").Value = "Armor"
$ws.Range("I33").Value = 1

$ws.Range("G34").Value = "Backup"
$ws.Range("I34").Value = 2

$ws.Range("G35").Value = "Cavalry"
$ws.Range("I35").Value = 3

$ws.Range("G36").Value = "Covert"
$ws.Range("I36").Value = 4

$ws.Range("G37").Value = "Dragon"
$ws.Range("I37").Value = 5

$ws.Range("G38").Value = "Flying"
$ws.Range("I38").Value = 6

$ws.Range("G39").Value = "Mystical"
$ws.Range("I39").Value = 7

$ws.Range("G40").Value = "QiAdept"
$ws.Range("I40").Value = 8

# ... "None" (value 0) added last, at the top of the ClassType block.
$ws.Range("G32").Value = "None"
$ws.Range("I32").Value = 0

# Reflect the final on-screen selection in the bottom (frozen) pane.
$ws.Range("H34").Select()
